$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new recording-session day (day 3, date 20230528) was added for subject 1.2.
# That pushes every row from the old row 7 onward down by one row, and the
# running index in column A gets renumbered sequentially (0..11).
#
# Rows 12 and 13 are brand new rows that did not exist before, so their
# column-A cell needs the bold/bordered/centered style ("A" style) copied
# over from an existing labelled cell before the value is written (Copy
# pulls in both value + format, so we overwrite the value afterwards).
$ws.Range("A11").Copy($ws.Range("A12"))
$ws.Range("A11").Copy($ws.Range("A13"))

# subject, date, -20_to_0, 0_to_10, date_int, cage  (col B..G); col A = index
$data = @(
    @(7,  5,  1.2, 20230528, 0.2377692307692308,     0.5196923076923078,  20230528, 1),
    @(8,  6,  1.3, 20230526, 0.002807692307692307,   0.000358974358974359,20230526, 1),
    @(9,  7,  1.3, 20230527, 0.120974358974359,      0.274025641025641,   20230527, 1),
    @(10, 8,  1.3, 20230528, 0.196051282051282,      0.3221025641025641,  20230528, 1),
    @(11, 9,  1.4, 20230526, 0.0002692307692307692,  0.000717948717948718,20230526, 1),
    @(12, 10, 1.4, 20230527, 0.1101025641025641,     0.2383589743589744,  20230527, 1),
    @(13, 11, 1.4, 20230528, 0.1194871794871795,     0.2931794871794872,  20230528, 1)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}
